$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume columns are treated as text, matching the original
# inline-string cell contents (avoids Excel auto-converting numeric-looking
# strings like "311.50" or "1.002" into numbers).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.061.45"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.59%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.811.16"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.50"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4618"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +4.72%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3761"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.62%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07404"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.84%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8636"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.58"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.813.98"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.642"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.384"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07077"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.64"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.64%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008726"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.87"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.076.00"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.324"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +3.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.040.58"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.17%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.86%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.263"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.85"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08911"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7730"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +4.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.172"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.70%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.517"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.901"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.42%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.126"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.64%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.76%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.66%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.233"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.36%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5288"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.73%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +13.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1674"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.647"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5022"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.96%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.31"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "104.60"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.45%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.669"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06332"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.19%  "

# Rows 27 and 28 swap coin identity (LidoDAOToken <-> EthereumClassic) with updated values
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.54"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.97%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.189"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.87%  "
